$wb = $excel.ActiveWorkbook

$wsSemi = $wb.Worksheets.Item("Semiconductors")

# Shift existing "Micron"/"MU" row (row 9) down to row 10
$micronName = $wsSemi.Range("B9").Value()
$micronTicker = $wsSemi.Range("C9").Value()
$wsSemi.Range("B10").Value = $micronName
$wsSemi.Range("C10").Value = $micronTicker

# Row 8 becomes Taiwan Semi with its local-listing ticker
$wsSemi.Range("B8").Value = "Taiwan Semi"
$wsSemi.Range("C8").Value = "2330 TT"

# Row 9 becomes Samsung with its local ticker
$wsSemi.Range("B9").Value = "Samsung"
$wsSemi.Range("C9").Value = "005930 KS"

# Row 11: new entry Broadcom / Avago
$wsSemi.Range("B11").Value = "Broadcom"
$wsSemi.Range("C11").Value = "Avago"

$wsSemi.Range("D11").Select()

$wsHw = $wb.Worksheets.Item("Hardware")
$wsHw.Activate()
